$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Ptprz1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.04741066666666666
$ws.Cells.Item(2, 8).Value = 0.142232
$ws.Cells.Item(2, 9).Value = 0.003188134523263584
$ws.Cells.Item(2, 10).Value = 0.003188134523263585
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.06163433333333333
$ws.Cells.Item(2, 14).Value = 0.184903
$ws.Cells.Item(2, 15).Value = 0.006690894379667537
$ws.Cells.Item(2, 16).Value = 0.006690894379667537
$ws.Cells.Item(2, 17).Value = 0.002922124832888888
$ws.Cells.Item(2, 18).Value = 0.026299123496
$ws.Cells.Item(2, 19).Value = 0.00002133147136332836
$ws.Cells.Item(2, 20).Value = 0.00002133147136332836

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Ptprz1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.04741066666666666
$ws.Cells.Item(3, 8).Value = 0.142232
$ws.Cells.Item(3, 9).Value = 0.003188134523263584
$ws.Cells.Item(3, 10).Value = 0.003188134523263585
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.1030763333333333
$ws.Cells.Item(3, 14).Value = 0.309229
$ws.Cells.Item(3, 15).Value = 0.01118975126488057
$ws.Cells.Item(3, 16).Value = 0.01118975126488057
$ws.Cells.Item(3, 17).Value = 0.004886917680888888
$ws.Cells.Item(3, 18).Value = 0.043982259128
$ws.Cells.Item(3, 19).Value = 0.00003567443231429812
$ws.Cells.Item(3, 20).Value = 0.00003567443231429812

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Ptprz1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.04741066666666666
$ws.Cells.Item(4, 8).Value = 0.142232
$ws.Cells.Item(4, 9).Value = 0.003188134523263584
$ws.Cells.Item(4, 10).Value = 0.003188134523263585
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 9.046962666666667
$ws.Cells.Item(4, 14).Value = 27.140888
$ws.Cells.Item(4, 15).Value = 0.9821193543554519
$ws.Cells.Item(4, 16).Value = 0.9821193543554518
$ws.Cells.Item(4, 17).Value = 0.4289225313351111
$ws.Cells.Item(4, 18).Value = 3.860302782016
$ws.Cells.Item(4, 19).Value = 0.003131128619585958
$ws.Cells.Item(4, 20).Value = 0.003131128619585958

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Ptprz1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.237305
$ws.Cells.Item(5, 8).Value = 0.711915
$ws.Cells.Item(5, 9).Value = 0.01595759596384214
$ws.Cells.Item(5, 10).Value = 0.01595759596384214
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.06163433333333333
$ws.Cells.Item(5, 14).Value = 0.184903
$ws.Cells.Item(5, 15).Value = 0.006690894379667537
$ws.Cells.Item(5, 16).Value = 0.006690894379667537
$ws.Cells.Item(5, 17).Value = 0.01462613547166666
$ws.Cells.Item(5, 18).Value = 0.131635219245
$ws.Cells.Item(5, 19).Value = 0.0001067705891474767
$ws.Cells.Item(5, 20).Value = 0.0001067705891474767

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Ptprz1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.237305
$ws.Cells.Item(6, 8).Value = 0.711915
$ws.Cells.Item(6, 9).Value = 0.01595759596384214
$ws.Cells.Item(6, 10).Value = 0.01595759596384214
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.1030763333333333
$ws.Cells.Item(6, 14).Value = 0.309229
$ws.Cells.Item(6, 15).Value = 0.01118975126488057
$ws.Cells.Item(6, 16).Value = 0.01118975126488057
$ws.Cells.Item(6, 17).Value = 0.02446052928166666
$ws.Cells.Item(6, 18).Value = 0.220144763535
$ws.Cells.Item(6, 19).Value = 0.0001785615296208557
$ws.Cells.Item(6, 20).Value = 0.0001785615296208557

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Ptprz1"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.237305
$ws.Cells.Item(7, 8).Value = 0.711915
$ws.Cells.Item(7, 9).Value = 0.01595759596384214
$ws.Cells.Item(7, 10).Value = 0.01595759596384214
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.046962666666667
$ws.Cells.Item(7, 14).Value = 27.140888
$ws.Cells.Item(7, 15).Value = 0.9821193543554519
$ws.Cells.Item(7, 16).Value = 0.9821193543554518
$ws.Cells.Item(7, 17).Value = 2.146889475613333
$ws.Cells.Item(7, 18).Value = 19.32200528052
$ws.Cells.Item(7, 19).Value = 0.0156722638450738
$ws.Cells.Item(7, 20).Value = 0.0156722638450738

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Ncam1"
$ws.Cells.Item(8, 3).Value = "Ptprz1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 14.58625866666667
$ws.Cells.Item(8, 8).Value = 43.758776
$ws.Cells.Item(8, 9).Value = 0.9808542695128942
$ws.Cells.Item(8, 10).Value = 0.9808542695128943
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.06163433333333333
$ws.Cells.Item(8, 14).Value = 0.184903
$ws.Cells.Item(8, 15).Value = 0.006690894379667537
$ws.Cells.Item(8, 16).Value = 0.006690894379667537
$ws.Cells.Item(8, 17).Value = 0.8990143287475554
$ws.Cells.Item(8, 18).Value = 8.091128958728
$ws.Cells.Item(8, 19).Value = 0.006562792319156732
$ws.Cells.Item(8, 20).Value = 0.006562792319156733

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Ncam1"
$ws.Cells.Item(9, 3).Value = "Ptprz1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 14.58625866666667
$ws.Cells.Item(9, 8).Value = 43.758776
$ws.Cells.Item(9, 9).Value = 0.9808542695128942
$ws.Cells.Item(9, 10).Value = 0.9808542695128943
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.1030763333333333
$ws.Cells.Item(9, 14).Value = 0.309229
$ws.Cells.Item(9, 15).Value = 0.01118975126488057
$ws.Cells.Item(9, 16).Value = 0.01118975126488057
$ws.Cells.Item(9, 17).Value = 1.503498060411555
$ws.Cells.Item(9, 18).Value = 13.531482543704
$ws.Cells.Item(9, 19).Value = 0.01097551530294542
$ws.Cells.Item(9, 20).Value = 0.01097551530294542

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Ncam1"
$ws.Cells.Item(10, 3).Value = "Ptprz1"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 14.58625866666667
$ws.Cells.Item(10, 8).Value = 43.758776
$ws.Cells.Item(10, 9).Value = 0.9808542695128942
$ws.Cells.Item(10, 10).Value = 0.9808542695128943
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 9.046962666666667
$ws.Cells.Item(10, 14).Value = 27.140888
$ws.Cells.Item(10, 15).Value = 0.9821193543554519
$ws.Cells.Item(10, 16).Value = 0.9821193543554518
$ws.Cells.Item(10, 17).Value = 131.9613376036764
$ws.Cells.Item(10, 18).Value = 1187.652038433088
$ws.Cells.Item(10, 19).Value = 0.9633159618907922
$ws.Cells.Item(10, 20).Value = 0.9633159618907922
